# Add a new "Rehab" / "Healthcare" keyword row into the Cluster_Keywords
# table on the Keywords sheet, just above the "Apart"/Housing block (i.e.
# as the new row 63), pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a fresh worksheet row at position 63 (shifts rows 63.. down by 1,
# and auto-extends the Cluster_Keywords table / named ranges / conditional
# formatting references that cover it).
$ws.Rows.Item(63).Insert()

# Populate the new row. Column B is the table's calculated "Length" column
# (LEN of the Stem), so set its formula the same way the other rows do.
$ws.Range("A63").Value = "Rehab"
$ws.Range("B63").Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"
$ws.Range("C63").Value = "Healthcare"

# Reset the view back to showing the top-left of the sheet (the source
# workbook had scrolled to A25 with A44 selected).
$win = $excel.Windows.Item(1)
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A1").Select()

$wb.Save()
